$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.298.85'
$ws.Range('E2').Value = '  -2.86%  '

$ws.Range('D3').Value = '1.733.02'
$ws.Range('E3').Value = '  -3.53%  '

$ws.Range('D4').Value = '''1.007'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  +0.12%  '

$ws.Range('D5').Value = '''323.25'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -4.51%  '

$ws.Range('D6').Value = '''1.003'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +0.14%  '

$ws.Range('D7').Value = '''0.4225'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -9.62%  '

$ws.Range('D8').Value = '''0.3598'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -3.46%  '

$ws.Range('D9').Value = '''44.97'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -1.32%  '

$ws.Range('D10').Value = '''0.07437'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -3.18%  '

$ws.Range('D11').Value = '''1.119'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -2.32%  '

$ws.Range('D12').Value = '''1.005'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +0.07%  '

$ws.Range('D13').Value = '''21.44'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -4.87%  '

$ws.Range('D14').Value = '''6.073'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -5.05%  '

$ws.Range('D15').Value = '''7.168'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -4.17%  '

$ws.Range('D16').Value = '1.735.54'
$ws.Range('E16').Value = '  -3.41%  '

$ws.Range('D17').Value = '''0.00001063'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -2.93%  '

$ws.Range('D18').Value = '''87.01'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +5.62%  '

$ws.Range('D19').Value = '''0.06019'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -10.85%  '

$ws.Range('D20').Value = '''1.002'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +0.08%  '

$ws.Range('D21').Value = '''16.79'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -4.11%  '

$ws.Range('D22').Value = '''6.057'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -6.03%  '

$ws.Range('D23').Value = '''0.5226'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -5.38%  '

$ws.Range('D24').Value = '27.365.87'
$ws.Range('E24').Value = '  -2.69%  '

$ws.Range('D25').Value = '''11.30'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -5.57%  '

$ws.Range('D26').Value = '''2.420'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +0.60%  '

$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').Value = '''20.08'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -3.74%  '

$ws.Range('B28').Value = 'LidoDAOToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D28').Value = '''2.372'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -1.45%  '

$ws.Range('D29').Value = '''149.58'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -1.15%  '

$ws.Range('D30').Value = '1.933.52'
$ws.Range('E30').Value = '  -3.54%  '

$ws.Range('D31').Value = '''1.270'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +0.63%  '

$ws.Range('D32').Value = '''126.37'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -5.94%  '

$ws.Range('D33').Value = '''3.738'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -7.86%  '

$ws.Range('D34').Value = '''5.584'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -6.07%  '

$ws.Range('D35').Value = '''0.09045'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -6.27%  '

$ws.Range('D36').Value = '''12.49'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +2.29%  '

$ws.Range('D37').Value = '''0.2147'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -3.75%  '

$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').Value = '''0.02278'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -5.04%  '

$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').Value = '''0.06139'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -4.01%  '

$ws.Range('B40').Value = 'InternetComputer(DFINITY)'
$ws.Range('C40').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D40').Value = '''5.012'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -5.02%  '

$ws.Range('B41').Value = 'TheSandbox'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D41').Value = '''0.6386'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -4.99%  '

$ws.Range('D42').Value = '''1.183'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -4.38%  '

$ws.Range('E43').Value = '  -4.86%  '

$ws.Range('D44').Value = '''1.002'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +0.03%  '

$ws.Range('D45').Value = '''7.865'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -2.67%  '

$ws.Range('D46').Value = '''13.56'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -4.92%  '

$ws.Range('D47').Value = '''3.738'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -3.28%  '

$ws.Range('D48').Value = '''0.5829'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -5.41%  '

$ws.Range('D49').Value = '''125.34'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -3.60%  '

$ws.Range('D50').Value = '''1.935'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -6.53%  '

$ws.Range('D51').Value = '''0.06840'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -4.03%  '
